$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "Updated cryptos list" refresh: new Price/Volume(1h) figures
# for each coin row, plus the NEARProtocol/Quant row order+data swap
# (rows 49 and 50). Values that look like plain numbers (e.g. "1.001")
# are written with a leading apostrophe so Excel keeps them as literal
# text, matching the original cells' text formatting, instead of
# auto-converting them to floating point numbers.
$ws.Range("D2").Value = "26.868.06"
$ws.Range("E2").Value = "  -1.83%  "
$ws.Range("D3").Value = "1.809.83"
$ws.Range("E3").Value = "  -0.87%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'309.90"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D7").Value = "'0.4642"
$ws.Range("E7").Value = "  +3.97%  "
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").Value = "'0.07355"
$ws.Range("E9").Value = "  -0.72%  "
$ws.Range("D10").Value = "'0.8778"
$ws.Range("E10").Value = "  +0.09%  "
$ws.Range("D11").Value = "'20.48"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "1.824.04"
$ws.Range("E12").Value = "  -0.35%  "
$ws.Range("D13").Value = "'5.361"
$ws.Range("E13").Value = "  -1.13%  "
$ws.Range("D14").Value = "'6.520"
$ws.Range("E14").Value = "  -2.78%  "
$ws.Range("D15").Value = "'91.87"
$ws.Range("E15").Value = "  -1.05%  "
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").Value = "'0.000008690"
$ws.Range("E18").Value = "  -1.45%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "'14.74"
$ws.Range("E20").Value = "  -2.20%  "
$ws.Range("D21").Value = "26.867.82"
$ws.Range("E21").Value = "  -1.88%  "
$ws.Range("D22").Value = "'5.316"
$ws.Range("E22").Value = "  -0.54%  "
$ws.Range("D23").Value = "'10.58"
$ws.Range("E23").Value = "  -3.22%  "
$ws.Range("D24").Value = "2.009.14"
$ws.Range("E24").Value = "  -2.39%  "
$ws.Range("D25").Value = "'1.897"
$ws.Range("E25").Value = "  -3.17%  "
$ws.Range("D26").Value = "'151.67"
$ws.Range("E26").Value = "  +0.37%  "
$ws.Range("D27").Value = "'18.40"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "'2.158"
$ws.Range("E28").Value = "  -4.99%  "
$ws.Range("D29").Value = "'5.329"
$ws.Range("E29").Value = "  -0.34%  "
$ws.Range("D30").Value = "'116.15"
$ws.Range("E30").Value = "  -0.90%  "
$ws.Range("D31").Value = "'0.08912"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").Value = "'0.7535"
$ws.Range("E32").Value = "  -5.12%  "
$ws.Range("D33").Value = "'1.159"
$ws.Range("E33").Value = "  -3.13%  "
$ws.Range("D34").Value = "'2.928"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "'4.466"
$ws.Range("E35").Value = "  -1.83%  "
$ws.Range("E36").Value = "  +0.14%  "
$ws.Range("D37").Value = "'1.103"
$ws.Range("E37").Value = "  -0.19%  "
$ws.Range("D38").Value = "'0.01966"
$ws.Range("E38").Value = "  -0.64%  "
$ws.Range("D39").Value = "'0.05251"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("D40").Value = "'2.419"
$ws.Range("E40").Value = "  +2.77%  "
$ws.Range("D41").Value = "'2.927"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").Value = "'0.5327"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "'7.175"
$ws.Range("E43").Value = "  -2.60%  "
$ws.Range("D44").Value = "'0.1665"
$ws.Range("E44").Value = "  -2.17%  "
$ws.Range("D45").Value = "'8.502"
$ws.Range("E45").Value = "  -2.03%  "
$ws.Range("D46").Value = "'0.4983"
$ws.Range("E46").Value = "  -2.17%  "
$ws.Range("D47").Value = "'10.31"
$ws.Range("E47").Value = "  -2.77%  "
$ws.Range("E48").Value = "  +0.14%  "
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D49").Value = "'1.670"
$ws.Range("E49").Value = "  -0.98%  "
$ws.Range("B50").Value = "Quant"
$ws.Range("C50").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D50").Value = "'103.75"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").Value = "'0.06295"
$ws.Range("E51").Value = "  -1.46%  "
